$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row values (row 1)
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Update row 2 values
$ws.Range("B2").Value = 16.403585443062106
$ws.Range("C2").Value = 41.524766107136323
$ws.Range("D2").Value = 16.859891899922577
$ws.Range("E2").Value = 31.536779922170894

# Update row 3 values
$ws.Range("B3").Value = 19.470153128188006
$ws.Range("C3").Value = 22.999451700416447
$ws.Range("D3").Value = 16.016676821730766
$ws.Range("E3").Value = 19.753652281622227

# Update selection to match new sqref B1:E3
$ws.Range("B1:E3").Select()
